$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenarios")
$ws.Range("C11").Value = ""
$ws.Range("C11").Interior.ThemeColor = 0
